# Update cryptocurrency price (D) and hourly volume-change (E) columns
# per the latest scrape, preserving each cell as literal text so
# Excel does not reinterpret numeric-looking price strings (e.g. "0.3410")
# or stomp the percentage strings padding/sign formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.399.15"
$ws.Range("E2").Value = "  -0.02%  "

$ws.Range("D3").Value = "1.573.63"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("E5").Value = "  +0.07%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.92"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3762"
$ws.Range("E7").Value = "  +2.44%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.88"
$ws.Range("E8").Value = "  +0.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3410"
$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.164"
$ws.Range("E10").Value = "  -0.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07666"
$ws.Range("E11").Value = "  +1.12%  "

$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.32"
$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.969"
$ws.Range("E14").Value = "  -1.83%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.914"
$ws.Range("E15").Value = "  +0.53%  "

$ws.Range("D16").Value = "1.573.49"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001136"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.54"
$ws.Range("E18").Value = "  +1.25%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06730"
$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("E20").Value = "  +0.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.73"
$ws.Range("E21").Value = "  +2.56%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.217"
$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5277"
$ws.Range("E23").Value = "  -4.74%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.98"
$ws.Range("E24").Value = "  +0.65%  "

$ws.Range("D25").Value = "22.405.63"
$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.404"
$ws.Range("E26").Value = "  -0.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.761"
$ws.Range("E27").Value = "  -6.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.26"
$ws.Range("E28").Value = "  +2.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "145.16"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.049"
$ws.Range("E30").Value = "  +2.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "126.03"
$ws.Range("E31").Value = "  +1.06%  "

$ws.Range("D32").Value = "1.746.89"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.200"
$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("E34").Value = "  +2.57%  "

$ws.Range("E35").Value = "  +4.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.04"
$ws.Range("E36").Value = "  -3.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08530"
$ws.Range("E37").Value = "  +0.58%  "

$ws.Range("E38").Value = "  +1.49%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2320"
$ws.Range("E39").Value = "  +1.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.352"
$ws.Range("E40").Value = "  +7.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06524"
$ws.Range("E41").Value = "  -0.05%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.472"
$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.61"
$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6479"
$ws.Range("E44").Value = "  +1.82%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.12"
$ws.Range("E45").Value = "  -2.69%  "

$ws.Range("E46").Value = "  +0.05%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6030"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.789"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.300"
$ws.Range("E49").Value = "  +9.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.097"
$ws.Range("E50").Value = "  -1.00%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "125.34"
$ws.Range("E51").Value = "  +3.74%  "
